# Updated symbol list on Sat Dec 17 03:51:37 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (D) and a handful of "Volume(1h)" (E) label cells
# with the latest scrape values. Rows 42/43 (CEJI / BKEXToken) also swap
# rank order, so their B/C/D/E cells are fully replaced together.
#
# All of these sheet cells are stored as text, so the number-formatted
# "Price" column must keep its values as plain text (not get auto-coerced
# into numeric values by Excel) -- force text format before assigning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, [string]$CellRef, [string]$Text)
    $rng = $Worksheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

# --- Simple Price (column D) refreshes -------------------------------
Set-TextValue $ws "D2"  "229.98"
Set-TextValue $ws "D3"  "22.37"
Set-TextValue $ws "D4"  "5.257"
Set-TextValue $ws "D5"  "0.05538"
Set-TextValue $ws "D6"  "3.381"
Set-TextValue $ws "D7"  "6.469"
Set-TextValue $ws "D9"  "0.7704"
Set-TextValue $ws "D10" "0.1377"
Set-TextValue $ws "D11" "0.07411"
Set-TextValue $ws "D12" "0.03151"
Set-TextValue $ws "D13" "0.02942"
Set-TextValue $ws "D14" "0.09260"
Set-TextValue $ws "D15" "0.001664"
Set-TextValue $ws "D16" "3.256"
Set-TextValue $ws "D17" "0.04781"
Set-TextValue $ws "D19" "0.006219"
Set-TextValue $ws "D20" "0.005231"
Set-TextValue $ws "D21" "0.001063"
Set-TextValue $ws "D22" "0.0001499"
Set-TextValue $ws "D23" "3.920"
Set-TextValue $ws "D26" "0.1243"
Set-TextValue $ws "D40" "0.03946"
Set-TextValue $ws "D41" "0.007122"
Set-TextValue $ws "D44" "0.008762"
Set-TextValue $ws "D45" "0.00005439"
Set-TextValue $ws "D46" "0.00000000750"
Set-TextValue $ws "D47" "0.7853"
Set-TextValue $ws "D49" "0.00002100"
Set-TextValue $ws "D50" "0.01010"

# --- Rows where both Price (D) and the Volume(1h) label (E) change ----
Set-TextValue $ws "D8"  "1.079"
Set-TextValue $ws "E8"  "7FTXTokenFTTWorstin24h"

Set-TextValue $ws "D18" "0.0005891"
Set-TextValue $ws "E18" "17OneONE"

Set-TextValue $ws "D27" "0.0005001"
Set-TextValue $ws "E27" "26UpBotsUBXTBestin24h"

Set-TextValue $ws "D48" "0.06675"
Set-TextValue $ws "E48" "47BOLOBOLO"

# --- Row 42 / 43: CEJI and BKEXToken swap rank/price ------------------
Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1035"
Set-TextValue $ws "E42" "41BKEXTokenBKK"

Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002708"
Set-TextValue $ws "E43" "42CEJICEJI"
